# Refresh the cryptos snapshot: updated prices / 1h volume figures,
# plus two rows (RenderToken/THORChain and VeChain/FTXToken) that
# swapped ranking order in the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-looking numbers (e.g. "229.65")
# that must stay TEXT, as in the source sheet (inline strings), not be
# coerced to numeric cells. Prefix with an apostrophe -- Excel's own
# 'force text' convention -- so COM stores it as text (apostrophe is
# not kept as part of the value).
function Set-TextPrice($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# Row 2
Set-TextPrice "D2" '41.401.78'
$ws.Range("E2").Value = '  +4.04%  '
# Row 3
Set-TextPrice "D3" '2.215.96'
$ws.Range("E3").Value = '  +2.49%  '
# Row 4
$ws.Range("E4").Value = '  +0.13%  '
# Row 5
Set-TextPrice "D5" '229.65'
# Row 6
Set-TextPrice "D6" '0.622'
$ws.Range("E6").Value = '  -0.81%  '
# Row 7
Set-TextPrice "D7" '60.97'
$ws.Range("E7").Value = '  -3.75%  '
# Row 8
$ws.Range("E8").Value = '  +0.05%  '
# Row 9
Set-TextPrice "D9" '0.399'
$ws.Range("E9").Value = '  +2.07%  '
# Row 10
Set-TextPrice "D10" '57.95'
$ws.Range("E10").Value = '  -1.77%  '
# Row 11
$ws.Range("E11").Value = '  +5.60%  '
# Row 12
$ws.Range("E12").Value = '  -0.36%  '
# Row 13
Set-TextPrice "D13" '2.545.44'
# Row 14
Set-TextPrice "D14" '15.55'
$ws.Range("E14").Value = '  -2.22%  '
# Row 15
$ws.Range("E15").Value = '  -1.39%  '
# Row 16
$ws.Range("E16").Value = '  -1.51%  '
# Row 17
$ws.Range("E17").Value = '  +0.50%  '
# Row 18
Set-TextPrice "D18" '2.212.59'
$ws.Range("E18").Value = '  +2.41%  '
# Row 19
Set-TextPrice "D19" '41.319.70'
$ws.Range("E19").Value = '  +4.13%  '
# Row 20
Set-TextPrice "D20" '72.52'
$ws.Range("E20").Value = '  +1.11%  '
# Row 21
Set-TextPrice "D21" '0.0₃0891'
$ws.Range("E21").Value = '  +5.54%  '
# Row 22
$ws.Range("E22").Value = '  +0.63%  '
# Row 23
Set-TextPrice "D23" '250.57'
$ws.Range("E23").Value = '  +9.10%  '
# Row 24
$ws.Range("E24").Value = '  +0.11%  '
# Row 25
Set-TextPrice "D25" '2.37'
$ws.Range("E25").Value = '  +0.35%  '
# Row 26
Set-TextPrice "D26" '2.26'
$ws.Range("E26").Value = '  -2.46%  '
# Row 27
$ws.Range("E27").Value = '  -0.15%  '
# Row 28
Set-TextPrice "D28" '167.19'
$ws.Range("E28").Value = '  -3.01%  '
# Row 29
$ws.Range("E29").Value = '  +0.68%  '
# Row 30
$ws.Range("E30").Value = '  +0.12%  '
# Row 31
$ws.Range("E31").Value = '  -2.03%  '
# Row 32
Set-TextPrice "D32" '2.55'
$ws.Range("E32").Value = '  -5.43%  '
# Row 33
$ws.Range("E33").Value = '  +0.01%  '
# Row 34
$ws.Range("E34").Value = '  +6.31%  '
# Row 35
$ws.Range("E35").Value = '  +0.70%  '
# Row 36
$ws.Range("E36").Value = '  +0.41%  '
# Row 37
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextPrice "D37" '6.53'
$ws.Range("E37").Value = '  -5.19%  '
# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextPrice "D38" '3.67'
$ws.Range("E38").Value = '  -0.48%  '
# Row 39
$ws.Range("E39").Value = '  -2.08%  '
# Row 40
$ws.Range("E40").Value = '  +0.25%  '
# Row 41
Set-TextPrice "D41" '0.000237'
$ws.Range("E41").Value = '  +27.28%  '
# Row 42
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextPrice "D42" '4.79'
$ws.Range("E42").Value = '  -5.55%  '
# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextPrice "D43" '0.0235'
$ws.Range("E43").Value = '  +3.88%  '
# Row 44
$ws.Range("E44").Value = '  +9.43%  '
# Row 45
Set-TextPrice "D45" '0.0976'
$ws.Range("E45").Value = '  +6.02%  '
# Row 46
Set-TextPrice "D46" '98.75'
$ws.Range("E46").Value = '  -3.63%  '
# Row 47
Set-TextPrice "D47" '1.20'
$ws.Range("E47").Value = '  -0.61%  '
# Row 48
Set-TextPrice "D48" '1.462.57'
$ws.Range("E48").Value = '  -3.48%  '
# Row 49
Set-TextPrice "D49" '16.50'
$ws.Range("E49").Value = '  -6.92%  '
# Row 50
$ws.Range("E50").Value = '  -0.86%  '
# Row 51
$ws.Range("E51").Value = '  -1.45%  '
